$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LECTURES")

# Insert a new column after the existing PLACE column (G) so that PLACE (G)
# becomes ROOMNUMBER and the newly inserted column (H) becomes BUILDINGNUMBER.
# LECTURER/GRADE_FACTOR shift right from H/I to I/J automatically.
$ws.Columns.Item(8).Insert()

# Rename headers.
$ws.Range("G1").Value = "ROOMNUMBER"
$ws.Range("H1").Value = "BUILDINGNUMBER"

# Replace the old textual PLACE codes with numeric room/building numbers.
# Rows 2-13 used to hold "MI01" -> room 102, building 5620.
$ws.Range("G2:G13").Value = 102
$ws.Range("H2:H13").Value = 5620

# Rows 14-22 used to hold "WI01" -> room 1801, building 5508.
$ws.Range("G14:G22").Value = 1801
$ws.Range("H14:H22").Value = 5508

# Rows 23-32 used to hold "ME01" -> room 2501, building 5101.
$ws.Range("G23:G32").Value = 2501
$ws.Range("H23:H32").Value = 5101
